$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-19 Wednesday" "2025-02-20 Thursday"

Replace-Text "653×7=" "446×2="
Replace-Text "950×4=" "780×7="
Replace-Text "473×9=" "133×9="
Replace-Text "696×7=" "799×7="
Replace-Text "774×5=" "641×3="
Replace-Text "515×2=" "298×4="
Replace-Text "631×5=" "524×7="
Replace-Text "482×7=" "364×3="
Replace-Text "465×3=" "304×7="
Replace-Text "118×7=" "592×4="
Replace-Text "219×9=" "791×4="
Replace-Text "890×4=" "945×5="
Replace-Text "981×8=" "348×8="
Replace-Text "188×8=" "994×9="
Replace-Text "610×7=" "587×7="
Replace-Text "792×6=" "949×2="
Replace-Text "651×2=" "855×3="
Replace-Text "733×6=" "310×8="
Replace-Text "148×4=" "825×5="
Replace-Text "114×7=" "804×6="
Replace-Text "556×4=" "892×8="
Replace-Text "233×9=" "948×8="
Replace-Text "833×5=" "996×5="
Replace-Text "753×7=" "103×7="
Replace-Text "221×7=" "945×7="

Write-Output "Done"
